$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 33 (older session entry): correct the error-message wording.
#    "out of VRAM" -> "error, out of memory"   (new shared string #206)
# ---------------------------------------------------------------------
$ws.Range("H33").Value = "error, out of memory"

# ---------------------------------------------------------------------
# 2) Row 37 (session 230910-0) gets its outcome filled in: it failed
#    with an out-of-memory error from slurm.
# ---------------------------------------------------------------------
$ws.Range("E37").Value = "failed, slurmstepd: out of memory"
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = "error, out of memory"
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = "NA"
$ws.Range("K37").Value = "NA"
$ws.Range("L37").Value = "NA"
$ws.Range("M37").Value = "NA"
$ws.Range("N37").Value = "NA"
$ws.Range("O37").Value = 0
$ws.Range("R37").Value = 6
$ws.Range("AB37").Value = "NA"
$ws.Range("AW37").Value = 1

# ---------------------------------------------------------------------
# 3) Insert a fresh row at 38 (pushing the stray note that used to sit
#    at row 39 down to row 44) to prepare the next session, 230910-1,
#    copying over the same "template" defaults row 37 used to carry.
# ---------------------------------------------------------------------
$ws.Range("A38:A42").EntireRow.Insert()
$ws.Rows("38:42").Clear()

$ws.Range("A38").Value = "230910-1"
$ws.Range("B38").Value = "train3dunet"
$ws.Range("C38").Value = "higher val eval score increase & better val prediction images by increasing train sample size by re-assigning the h5 files to train, val, and test sets (dataset04). Also, give more main memory in sbatch flags."
$ws.Range("D38").Value = "Better performance metrics & val prediction images, but not good enough to merit test3dunet session."
$ws.Range("E38").Value = "TBD"
$ws.Range("F38").Value = "TBD"
$ws.Range("G38").Value = "TBD"
$ws.Range("H38").Value = "TBD"
$ws.Range("I38").Value = "TBD"
$ws.Range("J38").Value = "TBD"
$ws.Range("K38").Value = "TBD"
$ws.Range("L38").Value = "TBD"
$ws.Range("M38").Value = "TBD"
$ws.Range("N38").Value = "TBD"
$ws.Range("O38").Value = "TBD"
$ws.Range("P38").Value = "dataset04"
$ws.Range("Q38").Value = 6
$ws.Range("R38").Value = "TBD"
$ws.Range("S38").Value = 5
$ws.Range("T38").Value = 1
$ws.Range("U38").Formula = "=S38 + T38"
$ws.Range("V38").Value = 1
$ws.Range("W38").Value = 1
$ws.Range("X38").Value = 16
$ws.Range("Y38").Value = "uint16"
$ws.Range("Z38").Value = 8
$ws.Range("AA38").Value = "uint8"
$ws.Range("AB38").Value = "TBD"
$ws.Range("AC38").Value = "NA"
$ws.Range("AD38").Value = "NA"
$ws.Range("AE38").Formula = "=1508.06553301511 + 0.00210606006752809 * (AM38*AN38*AO38) / 5 * U38"
$ws.Range("AF38").Value = 81920
$ws.Range("AG38").Value = 81049.6
$ws.Range("AH38").Value = 79.15
$ws.Range("AI38").Value = "NVIDIA A100-SXM4-80GB"
$ws.Range("AJ38").Value = 125
$ws.Range("AK38").Value = 1169
$ws.Range("AL38").Value = 414
$ws.Range("AM38").Value = 96
$ws.Range("AN38").Value = 784
$ws.Range("AO38").Value = 384
$ws.Range("AP38").Value = "yes"
$ws.Range("AQ38").Formula = "=_xlfn.FLOOR.MATH((AJ38 - AM38) / 2)"
$ws.Range("AR38").Formula = "=_xlfn.FLOOR.MATH((AK38 - AN38) / 2)"
$ws.Range("AS38").Formula = "=_xlfn.FLOOR.MATH((AL38 - AO38) / 2)"
$ws.Range("AT38").Value = "yes"
$ws.Range("AU38").Value = "patch = arbitrary even int_2^3 as large a possible"
$ws.Range("AV38").Value = "stride = floor (resolution - patch) / 2"
$ws.Range("AW38").Value = "TBD"
$ws.Range("AX38").Value = "TBD"
$ws.Range("AY38").Value = "TBD"

# ---------------------------------------------------------------------
# 4) Back to row 37: paste in the detailed slurm out-of-memory log last.
# ---------------------------------------------------------------------
$ws.Range("AX37").Value = "slurmstepd: error: Detected 2 oom-kill event(s) in StepId=5017888.batch. Some of your processes may have been killed by the cgroup out-of-memory handler."
$ws.Range("AY37").Value = "NA"

# ---------------------------------------------------------------------
# 5) Update the view so the newly added row is visible/selected.
# ---------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A38"), $true)
$ws.Range("H38").Select()
